$d = $word.ActiveDocument

# Paragraph 1: **ID__AFFARS_5315_topic_14__ID** -> **ID__AFFARS_5315_404_4__ID**
# Also drop the trailing space run that followed it.
$find = $d.Content.Find
$find.Execute("**ID__AFFARS_5315_topic_14__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5315_404_4__ID**", 2)

# Add a paragraph border (5 twips on each side) and widen the left indent.
$p1 = $d.Paragraphs.Item(1)
$p1.Format.LeftIndent = 11.25
$b = $p1.Format.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5
